$d = $word.ActiveDocument

# --- Locate the target paragraphs by scanning for the placeholder text "insert" ---
$conclusionParaIndex = -1
$yukiParaIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($conclusionParaIndex -eq -1 -and $t -eq "insert") {
        $conclusionParaIndex = $i
    }
    if ($yukiParaIndex -eq -1 -and $t -eq "Yuki did insert") {
        $yukiParaIndex = $i
    }
}

Write-Host "Conclusion placeholder paragraph index:" $conclusionParaIndex
Write-Host "Yuki placeholder paragraph index:" $yukiParaIndex

# --- Part 1: Replace the bold "insert" placeholder in the Conclusions/Discussion
#     section with the real first paragraph, then append four more BodyText
#     paragraphs with the rest of the conclusion/discussion text. ---
$p20 = $d.Paragraphs.Item($conclusionParaIndex)
$r = $p20.Range
$start = $r.Start
$end = $r.End - 1
$delRange = $d.Range($start, $end)
$delRange.Delete()
$insertPoint = $d.Range($start, $start)
$insertPoint.InsertAfter("In conclusion, our study aimed to predict Math, Reading, and Writing test scores based on a variety of personal and socio-economic variables. After attempting several modeling methods, we successfully identified several significant predictors for each test score. The variables: gender, ethnic group, parental education, lunch type, test preparation, parental marital status, sports participation, birth order, and weekly study hours were consistently found to influence Math and Reading Test Scores. The number of siblings and the mode of transportation variables were also significant predictors for Writing Test Scores along with the other predictors mentioned.")

$bodyTexts = @(
    "Our finalized models were determined using a combination of step-wise and criterion-based model selection methods which allowed us to understand the relationship between predictors and scores better. Our use of adjusted R-squared value as a criterion for model selection allows us to find a balance between bias and variance, ensuring that the models we choose are both robust and interpretable. We ran model diagnostics as well as influential observation diagnostics in order to confirm the validity of our finalized models. All of the results indicated the absence of worrisome influential observations.",
    "Our study findings revealed the complex relationship of variables that influence a student’s educational outcomes. For example, ethnic disparities, as evidenced by differences in test scores between different ethnic groups, are evident which reveal the need for targeted intervention in order to combat these inequalities that are seen in student’s educations.",
    "Some potential limitations of our study include negligence of interactions between variables and generalization. For this project, we specifically wanted to focus on implementing all the different types of model building techniques that we learned in P8130. Given the time constraint, we were unable to allocate more time to literature review to assess whether interactions between some of our covariates existed. In addition to this, we acknowledge that there was little to no background given for this dataset. We are unaware of the population that this specific dataset was sampled from so the results from this study are unable to be generalized to a larger population. Nevertheless, our rigorous modeling approaches and diagnostic checks do enhance the credibility of our findings.",
    "All in all, our study advances the general understanding of what predictors are important for Math, Reading and Writing test scores, providing a solid foundation for more future research and educational interventions."
)

$curIndex = $conclusionParaIndex
foreach ($txt in $bodyTexts) {
    $cp = $d.Paragraphs.Item($curIndex)
    $cp.Range.InsertParagraphAfter()
    $curIndex = $curIndex + 1
    $np = $d.Paragraphs.Item($curIndex)
    $np.Style = "BodyText"
    $np.Range.Text = $txt
}

# --- Part 2: Replace the three-run "Yuki did insert" paragraph with a single
#     plain run describing Yuki's actual contribution. Re-locate the paragraph by
#     its text instead of trusting a fixed offset, since four paragraphs were
#     inserted above and shifted all subsequent paragraph indices. ---
$yukiParaIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Yuki did insert") {
        $yukiParaIndex = $i
        break
    }
}
Write-Host "Re-located Yuki placeholder paragraph index:" $yukiParaIndex

$py = $d.Paragraphs.Item($yukiParaIndex)
Write-Host "Yuki paragraph text before replace:" $py.Range.Text
$ry = $py.Range
$ryStart = $ry.Start
$ryEnd = $ry.End - 1
$yukiDelRange = $d.Range($ryStart, $ryEnd)
$yukiDelRange.Delete()
$yukiInsertPoint = $d.Range($ryStart, $ryStart)
$yukiInsertPoint.InsertAfter("Yuki wrote the methods, majority of the results section, conclusion and discussion section and conducted data cleaning, exploration and visualization and SLR modeling.")

$pyAfter = $d.Paragraphs.Item($yukiParaIndex)
Write-Host "Yuki paragraph text after replace:" $pyAfter.Range.Text
Write-Host "Total paragraphs:" $d.Paragraphs.Count
